# Complete the second iteration of the project.
# Updates the "okokokok" attendance sheet: bumps the recorded date from
# 2021-11-11 to 2021-11-12 for both tracked people, and fills in the
# "Answered" / "QA" counters for each of the four Date/Answered/QA
# triplets (columns C-E, F-H, I-K, L-N) on rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("okokokok")

# Helper-free approach: every literal we write (dates like "2021-11-12",
# and small numeric-looking tallies like "1" / "0") would otherwise be
# auto-converted by the typed-input parser (dates -> serial numbers,
# digit strings -> numbers) and would stamp the cell with a brand-new
# number-format style. Routing the literal through a text Formula and
# then collapsing it to a static value with PasteSpecial(xlPasteValues)
# keeps it as plain text and leaves the cell's style untouched (style 0).

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $escaped = $text -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy()
    $rng.PasteSpecial(-4163)
}

# Row 2 (raghu_palash / Palash)
Set-TextValue "C2" "2021-11-12"
Set-TextValue "D2" "1, 2"
Set-TextValue "E2" "2"
Set-TextValue "F2" "2021-11-12"
Set-TextValue "G2" "1"
Set-TextValue "H2" "1"
Set-TextValue "I2" "2021-11-12"
Set-TextValue "J2" "1"
Set-TextValue "K2" "1"
Set-TextValue "L2" "2021-11-12"
Set-TextValue "M2" "1, 2"
Set-TextValue "N2" "2"

# Row 3 (lalsa04 / लालसा😀)
Set-TextValue "C3" "2021-11-12"
Set-TextValue "E3" "0"
Set-TextValue "F3" "2021-11-12"
Set-TextValue "G3" "2, 3"
Set-TextValue "H3" "2"
Set-TextValue "I3" "2021-11-12"
Set-TextValue "K3" "0"
Set-TextValue "L3" "2021-11-12"
Set-TextValue "M3" "1"
Set-TextValue "N3" "1"
